$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string (row 1, col A)
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 10:25"

# 2. Update country statistics rows that changed in place

# India (row 5)
$ws.Range("B5").Value = 7708947
$ws.Range("C5").Value = 3789
$ws.Range("D5").Value = 6874518
$ws.Range("E5").Value = 717748
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 28
$ws.Range("H5").Value = 116681

# Rusia (row 7)
$ws.Range("B7").Value = 1463306
$ws.Range("C7").Value = 15971
$ws.Range("D7").Value = 1107988
$ws.Range("E7").Value = 330076
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 290
$ws.Range("H7").Value = 25242

# Filipinas (row 23)
$ws.Range("B23").Value = 363888
$ws.Range("C23").Value = 1664
$ws.Range("D23").Value = 312333
$ws.Range("E23").Value = 44772
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 38
$ws.Range("H23").Value = 6783

# Canada (row 32)
$ws.Range("B32").Value = 206360
$ws.Range("C32").Value = 406
$ws.Range("D32").Value = 173748
$ws.Range("E32").Value = 22783
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 9829

# Polonia (row 33)
$ws.Range("B33").Value = 202579
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 102204
$ws.Range("E33").Value = 96524
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 3851

# Singapur (row 65)
$ws.Range("B65").Value = 57941
$ws.Range("C65").Value = 8
$ws.Range("D65").Value = 57821
$ws.Range("E65").Value = 92
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 28

# Afganistan (row 78)
$ws.Range("B78").Value = 40626
$ws.Range("C78").Value = 116
$ws.Range("D78").Value = 33831
$ws.Range("E78").Value = 5290
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 1505

# Eslovaquia (row 83)
$ws.Range("B83").Value = 35330
$ws.Range("C83").Value = 1728
$ws.Range("D83").Value = 8763
$ws.Range("E83").Value = 26452
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 17
$ws.Range("H83").Value = 115

# 3. Countries list reorder: Estonia moves up, now sits between Aruba and Islandia.
#    Old order (rows 141-145): Aruba, Islandia, Mayotte, Estonia, Somalia
#    New order (rows 141-145): Aruba, Estonia, Islandia, Mayotte, Somalia
#    Row 141 (Aruba) and row 145 (Somalia) are unchanged.
#    Islandia's old stats move to row 143, Mayotte's old stats move to row 144,
#    and row 142 gets Estonia's updated stats.

# Row 142 becomes Estonia (new/updated totals)
$ws.Range("A142").Value = "Estonia"
$ws.Range("B142").Value = 4247
$ws.Range("C142").Value = 76
$ws.Range("D142").Value = 3366
$ws.Range("E142").Value = 810
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 71

# Row 143 becomes Islandia (previously row 142's totals, unchanged values)
$ws.Range("A143").Value = "Islandia"
$ws.Range("B143").Value = 4230
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 3013
$ws.Range("E143").Value = 1206
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 11

# Row 144 becomes Mayotte (previously row 143's totals, unchanged values)
$ws.Range("A144").Value = "Mayotte"
$ws.Range("B144").Value = 4203
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 2964
$ws.Range("E144").Value = 1195
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 44

# Taiwan (row 178)
$ws.Range("B178").Value = 548
$ws.Range("C178").Value = 4
$ws.Range("D178").Value = 497
$ws.Range("E178").Value = 44
$ws.Range("F178").Value = 0

# Brunei (row 196)
$ws.Range("B196").Value = 148
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 143
$ws.Range("E196").Value = 2
$ws.Range("F196").Value = 0
